# Apply updated dSF (column F) values to the active worksheet.
# These are the "repulled" / recalculated values for the dSF column
# (rows 2-28) as described in the commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = -1
    3  = 4
    4  = 2
    5  = 7
    6  = -3
    7  = 6
    8  = -1
    9  = 6
    10 = 1
    11 = -3
    12 = 0
    15 = 2
    18 = -1
    19 = 2
    20 = -1
    21 = -1
    22 = -2
    23 = -2
    24 = 1
    25 = 2
    26 = -2
    27 = -3
    28 = -2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
